# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# OFF sheet ("R" row, row 3): update Short Att / Short Comp / Deep Att / Deep Comp
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 304
$wsOff.Range("C3").Value = 194
$wsOff.Range("D3").Value = 143
$wsOff.Range("E3").Value = 61

# DEF sheet ("R" row, row 3): update Short Att / Short Comp / Deep Att / Deep Comp / Short Int / Deep Int
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 395
$wsDef.Range("C3").Value = 287
$wsDef.Range("D3").Value = 115
$wsDef.Range("E3").Value = 54
$wsDef.Range("F3").Value = 7
$wsDef.Range("G3").Value = 2
